$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Donor cell (untouched, style s="3") used to restore formatting after
# writing percentage strings, since Excel auto-parses "NN%" into a number
# with a new percent number format otherwise.
$fmtDonor = $ws.Cells.Item(5, 8)

$ws.Cells.Item(2, 5).Value = '2026-02-04 18:50:25'
$ws.Cells.Item(2, 15).Value = '-1.1 °C'
$ws.Cells.Item(3, 5).Value = '2026-02-04 18:50:28'
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = '89%'
$fmtDonor.Copy()
$ws.Cells.Item(3, 8).PasteSpecial(-4122)
$ws.Cells.Item(3, 15).Value = '-4.7 °C'
$ws.Cells.Item(4, 5).Value = '2026-02-04 18:50:30'
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = '81%'
$fmtDonor.Copy()
$ws.Cells.Item(4, 8).PasteSpecial(-4122)
$ws.Cells.Item(4, 10).Value = '993.2 hPa'
$ws.Cells.Item(4, 15).Value = '6.1 °C'
$ws.Cells.Item(5, 5).Value = '2026-02-04 18:50:33'
$ws.Cells.Item(6, 5).Value = '2026-02-04 18:50:36'
$ws.Cells.Item(6, 15).Value = '11.1 °C'
$ws.Cells.Item(7, 5).Value = '2026-02-04 18:50:38'
$ws.Cells.Item(8, 5).Value = '2026-02-04 18:50:41'
$ws.Cells.Item(8, 15).Value = '7.7 °C'
$ws.Cells.Item(9, 5).Value = '2026-02-04 18:50:44'
$ws.Cells.Item(10, 5).Value = '2026-02-04 18:50:46'
$ws.Cells.Item(11, 5).Value = '2026-02-04 18:50:49'
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = '82%'
$fmtDonor.Copy()
$ws.Cells.Item(11, 8).PasteSpecial(-4122)
$ws.Cells.Item(11, 15).Value = '0.8 °C'
$ws.Cells.Item(12, 5).Value = '2026-02-04 18:50:52'
$ws.Cells.Item(13, 5).Value = '2026-02-04 18:50:55'
$ws.Cells.Item(14, 5).Value = '2026-02-04 18:50:57'
$ws.Cells.Item(14, 13).Value = '-2.7 °C 18:27 TU'
$ws.Cells.Item(14, 15).Value = '-6.2 °C'
$ws.Cells.Item(15, 5).Value = '2026-02-04 18:51:00'
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = '81%'
$fmtDonor.Copy()
$ws.Cells.Item(15, 8).PasteSpecial(-4122)
$ws.Cells.Item(16, 5).Value = '2026-02-04 18:51:02'
$ws.Cells.Item(17, 5).Value = '2026-02-04 18:51:05'
$ws.Cells.Item(18, 5).Value = '2026-02-04 18:51:08'
$ws.Cells.Item(19, 5).Value = '2026-02-04 18:51:10'
$ws.Cells.Item(20, 5).Value = '2026-02-04 18:51:13'
$ws.Cells.Item(20, 8).NumberFormat = "@"
$ws.Cells.Item(20, 8).Value = '84%'
$fmtDonor.Copy()
$ws.Cells.Item(20, 8).PasteSpecial(-4122)
$ws.Cells.Item(21, 5).Value = '2026-02-04 18:51:16'
$ws.Cells.Item(21, 10).Value = '992.9 hPa'
$ws.Cells.Item(22, 5).Value = '2026-02-04 18:51:19'
$ws.Cells.Item(23, 5).Value = '2026-02-04 18:51:21'
$ws.Cells.Item(23, 10).Value = '992.3 hPa'
$ws.Cells.Item(24, 5).Value = '2026-02-04 18:51:24'
$ws.Cells.Item(24, 10).Value = '991.4 hPa'
$ws.Cells.Item(24, 15).Value = '10.1 °C'
$ws.Cells.Item(25, 5).Value = '2026-02-04 18:51:27'
$ws.Cells.Item(26, 5).Value = '2026-02-04 18:51:29'
$ws.Cells.Item(27, 5).Value = '2026-02-04 18:51:32'
$ws.Cells.Item(27, 10).Value = '992.9 hPa'
$ws.Cells.Item(27, 15).Value = '10.7 °C'
$ws.Cells.Item(28, 5).Value = '2026-02-04 18:51:34'
$ws.Cells.Item(29, 5).Value = '2026-02-04 18:51:37'
$ws.Cells.Item(30, 5).Value = '2026-02-04 18:51:40'
$ws.Cells.Item(30, 15).Value = '-5.3 °C'
$ws.Cells.Item(31, 5).Value = '2026-02-04 18:51:42'
$ws.Cells.Item(31, 15).Value = '4.3 °C'
$ws.Cells.Item(32, 5).Value = '2026-02-04 18:51:45'
$ws.Cells.Item(32, 10).Value = '993.5 hPa'
$ws.Cells.Item(33, 5).Value = '2026-02-04 18:51:48'
$ws.Cells.Item(34, 5).Value = '2026-02-04 18:51:50'
$ws.Cells.Item(34, 15).Value = '3.6 °C'
$ws.Cells.Item(35, 5).Value = '2026-02-04 18:51:53'
$ws.Cells.Item(36, 5).Value = '2026-02-04 18:51:56'

$excel.CutCopyMode = 0
